$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update sheet dimension happens automatically when cells are written ---

# --- Headers: I1 = "I0", J1 = "IF" (copy header style from H1) ---
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-77: columns I (I0) and J (IF) ---
$arr = New-Object 'object[,]' 76,2
$arr[0,0] = 8; $arr[0,1] = 8
$arr[1,0] = 8; $arr[1,1] = 8
$arr[2,0] = 5; $arr[2,1] = 6
$arr[3,0] = 7; $arr[3,1] = 7
$arr[4,0] = 9; $arr[4,1] = 9
$arr[5,0] = 7; $arr[5,1] = 7
$arr[6,0] = 8; $arr[6,1] = 8
$arr[7,0] = 9; $arr[7,1] = 9
$arr[8,0] = 7; $arr[8,1] = 7
$arr[9,0] = 6; $arr[9,1] = 7
$arr[10,0] = 8; $arr[10,1] = 8
$arr[11,0] = 7; $arr[11,1] = 7
$arr[12,0] = 8; $arr[12,1] = 8
$arr[13,0] = 8; $arr[13,1] = 8
$arr[14,0] = 7; $arr[14,1] = 7
$arr[15,0] = 8; $arr[15,1] = 8
$arr[16,0] = 8; $arr[16,1] = 8
$arr[17,0] = 8; $arr[17,1] = 8
$arr[18,0] = 7; $arr[18,1] = 7
$arr[19,0] = 8; $arr[19,1] = 8
$arr[20,0] = 6; $arr[20,1] = 6
$arr[21,0] = 7; $arr[21,1] = 7
$arr[22,0] = 6; $arr[22,1] = 6
$arr[23,0] = 7; $arr[23,1] = 8
$arr[24,0] = 9; $arr[24,1] = 9
$arr[25,0] = 8; $arr[25,1] = 8
$arr[26,0] = 7; $arr[26,1] = 7
$arr[27,0] = 10; $arr[27,1] = 10
$arr[28,0] = 7; $arr[28,1] = 7
$arr[29,0] = 8; $arr[29,1] = 8
$arr[30,0] = 7; $arr[30,1] = 7
$arr[31,0] = 6; $arr[31,1] = 6
$arr[32,0] = 7; $arr[32,1] = 7
$arr[33,0] = 6; $arr[33,1] = 6
$arr[34,0] = 7; $arr[34,1] = 7
$arr[35,0] = 7; $arr[35,1] = 7
$arr[36,0] = 8; $arr[36,1] = 8
$arr[37,0] = 7; $arr[37,1] = 7
$arr[38,0] = 10; $arr[38,1] = 10
$arr[39,0] = 6; $arr[39,1] = 6
$arr[40,0] = 9; $arr[40,1] = 9
$arr[41,0] = 8; $arr[41,1] = 8
$arr[42,0] = 8; $arr[42,1] = 8
$arr[43,0] = 10; $arr[43,1] = 10
$arr[44,0] = 7; $arr[44,1] = 7
$arr[45,0] = 6; $arr[45,1] = 6
$arr[46,0] = 7; $arr[46,1] = 7
$arr[47,0] = 8; $arr[47,1] = 8
$arr[48,0] = 7; $arr[48,1] = 7
$arr[49,0] = 8; $arr[49,1] = 8
$arr[50,0] = 8; $arr[50,1] = 8
$arr[51,0] = 8; $arr[51,1] = 8
$arr[52,0] = 7; $arr[52,1] = 7
$arr[53,0] = 8; $arr[53,1] = 8
$arr[54,0] = 7; $arr[54,1] = 7
$arr[55,0] = 9; $arr[55,1] = 9
$arr[56,0] = 5; $arr[56,1] = 6
$arr[57,0] = 9; $arr[57,1] = 9
$arr[58,0] = 9; $arr[58,1] = 9
$arr[59,0] = 7; $arr[59,1] = 8
$arr[60,0] = 9; $arr[60,1] = 9
$arr[61,0] = 9; $arr[61,1] = 9
$arr[62,0] = 8; $arr[62,1] = 9
$arr[63,0] = 6; $arr[63,1] = 6
$arr[64,0] = 8; $arr[64,1] = 8
$arr[65,0] = 4; $arr[65,1] = 6
$arr[66,0] = 9; $arr[66,1] = 9
$arr[67,0] = 9; $arr[67,1] = 9
$arr[68,0] = 8; $arr[68,1] = 8
$arr[69,0] = 8; $arr[69,1] = 8
$arr[70,0] = 7; $arr[70,1] = 7
$arr[71,0] = 8; $arr[71,1] = 8
$arr[72,0] = 4; $arr[72,1] = 5
$arr[73,0] = 3; $arr[73,1] = 4
$arr[74,0] = 4; $arr[74,1] = 4
$arr[75,0] = 3; $arr[75,1] = 4

$ws.Range("I2:J77").Value = $arr
